$wb = $excel.ActiveWorkbook

# ALC sheet
$ws = $wb.Worksheets.Item(1)
$ws.Range("H5").Value = 3000167.2
$ws.Range("I5").Value = 4500250
$ws.Range("K5").Value = 4500250
$ws.Range("M5").Value = -4500135
$ws.Range("H28").Value = 1333.0952
$ws.Range("J28").Value = 3092
$ws.Range("L28").Value = 3092
$ws.Range("N28").Value = -4062
$ws.Range("H62").Value = 1800
$ws.Range("J62").Value = 2500
$ws.Range("L62").Value = 2500
$ws.Range("N62").Value = -3748
$ws.Range("H64").Value = 5337.25
$ws.Range("I64").Value = 5448.8335
$ws.Range("K64").Value = 5448.8335
$ws.Range("M64").Value = -5200.8335
$ws.Range("H65").Value = 1800
$ws.Range("J65").Value = 2500
$ws.Range("L65").Value = 12500
$ws.Range("N65").Value = -18740
$ws.Range("H67").Value = 5337.25
$ws.Range("I67").Value = 5448.8335
$ws.Range("K67").Value = 5448.8335
$ws.Range("M67").Value = -4590.8335
$ws.Range("H76").Value = 14687.8
$ws.Range("I76").Value = 14283.143
$ws.Range("K76").Value = 14283.143
$ws.Range("M76").Value = -13968.143
$ws.Range("H79").Value = 14687.8
$ws.Range("I79").Value = 14283.143
$ws.Range("K79").Value = 14283.143
$ws.Range("M79").Value = -13191.143
$ws.Range("H96").Value = 1294.9
$ws.Range("J96").Value = 3014.5
$ws.Range("L96").Value = 9043.5
$ws.Range("N96").Value = -11789.5
$ws.Range("H107").Value = 1297.1333
$ws.Range("I107").Value = 1114.3636
$ws.Range("K107").Value = 1114.3636
$ws.Range("M107").Value = 805.6364000000001
$ws.Range("H112").Value = 4352.7144
$ws.Range("I112").Value = 989.6667
$ws.Range("J112").Value = 6875
$ws.Range("K112").Value = 2969.0001
$ws.Range("L112").Value = 20625
$ws.Range("M112").Value = -1861.0001
$ws.Range("N112").Value = -22841
$ws.Range("H115").Value = 9949.5
$ws.Range("I115").Value = 9899.5
$ws.Range("J115").Value = 9999.5
$ws.Range("K115").Value = 29698.5
$ws.Range("L115").Value = 29998.5
$ws.Range("M115").Value = -28131.5
$ws.Range("N115").Value = -33132.5
$ws.Range("H117").Value = 40000
$ws.Range("I117").Value = 40000
$ws.Range("K117").Value = 40000
$ws.Range("M117").Value = -35411
$ws.Range("H125").Value = 1074.75
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920
$ws.Range("H138").Value = 2159.1794
$ws.Range("I138").Value = 1310.9445
$ws.Range("K138").Value = 3932.8335
$ws.Range("M138").Value = 1207.1665
$ws.Range("H140").Value = 114632.336
$ws.Range("J140").Value = 114632.336
$ws.Range("L140").Value = 114632.336
$ws.Range("N140").Value = -124992.336

# ARM sheet
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 3176.342
$ws.Range("I61").Value = 2920.6128
$ws.Range("K61").Value = 2920.6128
$ws.Range("M61").Value = -2708.6128
$ws.Range("H102").Value = 9312.700000000001
$ws.Range("I102").Value = 9312.700000000001
$ws.Range("K102").Value = 9312.700000000001
$ws.Range("M102").Value = -7690.700000000001
$ws.Range("H132").Value = 2388.6206
$ws.Range("I132").Value = 2368
$ws.Range("K132").Value = 7104
$ws.Range("M132").Value = -4574
$ws.Range("H134").Value = 96666.664
$ws.Range("J134").Value = 96666.664
$ws.Range("L134").Value = 96666.664
$ws.Range("N134").Value = -106806.664
$ws.Range("H136").Value = 3176.342
$ws.Range("I136").Value = 2920.6128
$ws.Range("K136").Value = 8761.838400000001
$ws.Range("M136").Value = -6211.838400000001
$ws.Range("H140").Value = 97990
$ws.Range("J140").Value = 97990
$ws.Range("L140").Value = 97990
$ws.Range("N140").Value = -108350

# BSM sheet
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 2220.4075
$ws.Range("I20").Value = 1972.2
$ws.Range("K20").Value = 1972.2
$ws.Range("M20").Value = -1725.2
$ws.Range("H99").Value = 3106.0667
$ws.Range("I99").Value = 2613.7856
$ws.Range("J99").Value = 9998
$ws.Range("K99").Value = 2613.7856
$ws.Range("L99").Value = 9998
$ws.Range("M99").Value = -1115.7856
$ws.Range("N99").Value = -12994
$ws.Range("H107").Value = 2797.6667
$ws.Range("I107").Value = 2658.1428
$ws.Range("K107").Value = 2658.1428
$ws.Range("M107").Value = -738.1428000000001

# CRP sheet
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 35487
$ws.Range("I16").Value = 649.6667
$ws.Range("J16").Value = 139999
$ws.Range("K16").Value = 649.6667
$ws.Range("L16").Value = 139999
$ws.Range("M16").Value = -362.6667
$ws.Range("N16").Value = -140573
$ws.Range("H94").Value = 1466.6666
$ws.Range("I94").Value = 1350
$ws.Range("K94").Value = 1350
$ws.Range("M94").Value = -899
$ws.Range("H112").Value = 100000
$ws.Range("J112").Value = 100000
$ws.Range("L112").Value = 100000
$ws.Range("N112").Value = -102954
$ws.Range("H113").Value = 35487
$ws.Range("I113").Value = 649.6667
$ws.Range("J113").Value = 139999
$ws.Range("K113").Value = 649.6667
$ws.Range("L113").Value = 139999
$ws.Range("M113").Value = 1520.3333
$ws.Range("N113").Value = -144339
$ws.Range("H138").Value = 94998.5
$ws.Range("J138").Value = 94998.5
$ws.Range("L138").Value = 94998.5
$ws.Range("N138").Value = -105278.5

# CUL sheet
$ws = $wb.Worksheets.Item(5)
$ws.Range("H10").Value = 445.33334
$ws.Range("I10").Value = 231.8
$ws.Range("J10").Value = 712.25
$ws.Range("K10").Value = 695.4000000000001
$ws.Range("L10").Value = 2136.75
$ws.Range("M10").Value = -556.4000000000001
$ws.Range("N10").Value = -2414.75
$ws.Range("H106").Value = 19833.334
$ws.Range("J106").Value = 19833.334
$ws.Range("L106").Value = 59500.00199999999
$ws.Range("N106").Value = -61392.00199999999

# GSM sheet
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 3869
$ws.Range("I80").Value = 4062.8
$ws.Range("K80").Value = 4062.8
$ws.Range("M80").Value = -3064.8
$ws.Range("H83").Value = 3869
$ws.Range("I83").Value = 4062.8
$ws.Range("K83").Value = 20314
$ws.Range("M83").Value = -15322
$ws.Range("H99").Value = 24981.666
$ws.Range("I99").Value = 9985.25
$ws.Range("K99").Value = 9985.25
$ws.Range("M99").Value = -7739.25
$ws.Range("H113").Value = 7287
$ws.Range("J113").Value = 6233.3335
$ws.Range("L113").Value = 6233.3335
$ws.Range("N113").Value = -10573.3335
$ws.Range("H135").Value = 294187.84
$ws.Range("J135").Value = 87376.664
$ws.Range("L135").Value = 87376.664
$ws.Range("N135").Value = -97516.664

# LTW sheet
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1215.35
$ws.Range("J16").Value = 1600.3334
$ws.Range("L16").Value = 1600.3334
$ws.Range("N16").Value = -1940.3334
$ws.Range("H40").Value = 2848.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2848.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2848.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3120.5
$ws.Range("H46").Value = 5303.6113
$ws.Range("I46").Value = 6849.909
$ws.Range("J46").Value = 2873.7144
$ws.Range("K46").Value = 6849.909
$ws.Range("L46").Value = 2873.7144
$ws.Range("M46").Value = -6661.909
$ws.Range("N46").Value = -3249.7144
$ws.Range("H55").Value = 425.91177
$ws.Range("I55").Value = 510.53845
$ws.Range("K55").Value = 510.53845
$ws.Range("M55").Value = -337.53845
$ws.Range("H122").Value = 9966.375
$ws.Range("J122").Value = 3247
$ws.Range("L122").Value = 9741
$ws.Range("N122").Value = -14641
$ws.Range("H125").Value = 89999
$ws.Range("J125").Value = 89999
$ws.Range("L125").Value = 89999
$ws.Range("N125").Value = -99839
$ws.Range("H132").Value = 2467.342
$ws.Range("I132").Value = 2245.52
$ws.Range("J132").Value = 2893.923
$ws.Range("K132").Value = 6736.559999999999
$ws.Range("L132").Value = 8681.769
$ws.Range("M132").Value = -4206.559999999999
$ws.Range("N132").Value = -13741.769
$ws.Range("H141").Value = 88297.7
$ws.Range("J141").Value = 88664.11
$ws.Range("L141").Value = 88664.11
$ws.Range("N141").Value = -99024.11

# WVR sheet
$ws = $wb.Worksheets.Item(8)
$ws.Range("H141").Value = 56757
$ws.Range("J141").Value = 59106.332
$ws.Range("L141").Value = 59106.332
$ws.Range("N141").Value = -69466.33199999999

